$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (score) and E (reviews count) hold numeric-looking text values
# (e.g. "8.3", "6,609") that must stay text, matching the source data, so
# force a text format on each such cell before assigning its new value.
foreach ($addr in @("C2", "E2", "C3", "E3", "C4", "E4", "C5", "E5", "C6", "E6", "C7", "E7", "E8", "C9", "E9", "C10", "E10", "C11", "E11", "C12", "E12", "C13", "E13", "C14", "E14", "C15", "E15", "C16", "E16", "C17", "E17", "C18", "E18", "C19", "E19", "C20", "E20", "C21", "E21", "C22", "E22", "C23", "E23", "C24", "E24", "C25", "E25", "C26", "E26")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = 'The People - Paris Marais'
$ws.Range("B2").Value = 'US$392'
$ws.Range("C2").Value = '8.3'
$ws.Range("D2").Value = 'Very Good'
$ws.Range("E2").Value = '6,609'

$ws.Range("A3").Value = 'The People - Paris Nation'
$ws.Range("B3").Value = 'US$275'
$ws.Range("C3").Value = '8.4'
$ws.Range("D3").Value = 'Very Good'
$ws.Range("E3").Value = '6,122'

$ws.Range("A4").Value = 'Hotel The Playce by Happyculture'
$ws.Range("B4").Value = 'US$1,009'
$ws.Range("C4").Value = '7.7'
$ws.Range("D4").Value = 'Good'
$ws.Range("E4").Value = '2,242'

$ws.Range("A5").Value = 'Hôtel des Andelys'
$ws.Range("B5").Value = 'US$514'
$ws.Range("C5").Value = '4.0'
$ws.Range("D5").Value = 'Review score'
$ws.Range("E5").Value = '928'

$ws.Range("A6").Value = 'The People Paris Belleville'
$ws.Range("B6").Value = 'US$282'
$ws.Range("C6").Value = '8.2'
$ws.Range("D6").Value = 'Very Good'
$ws.Range("E6").Value = '2,735'

$ws.Range("A7").Value = 'Austin''s Saint Lazare Hotel'
$ws.Range("B7").Value = 'US$1,489'
$ws.Range("C7").Value = '8.1'
$ws.Range("D7").Value = 'Very Good'
$ws.Range("E7").Value = '2,537'

$ws.Range("A8").Value = 'Hôtel Le Daum'
$ws.Range("B8").Value = 'US$1,534'
$ws.Range("E8").Value = '943'

$ws.Range("A9").Value = 'Hôtel La Conversation'
$ws.Range("B9").Value = 'US$1,323'
$ws.Range("C9").Value = '8.2'
$ws.Range("D9").Value = 'Very Good'
$ws.Range("E9").Value = '949'

$ws.Range("A10").Value = 'Hotel Royal Phare'
$ws.Range("B10").Value = 'US$1,506'
$ws.Range("C10").Value = '8.3'
$ws.Range("D10").Value = 'Very Good'
$ws.Range("E10").Value = '1,734'

$ws.Range("A11").Value = 'Hôtel de l''Europe'
$ws.Range("B11").Value = 'US$605'
$ws.Range("C11").Value = '5.3'
$ws.Range("D11").Value = 'Review score'
$ws.Range("E11").Value = '476'

$ws.Range("A12").Value = 'Enjoy Hostel'
$ws.Range("B12").Value = 'US$306'
$ws.Range("C12").Value = '6.9'
$ws.Range("D12").Value = 'Review score'
$ws.Range("E12").Value = '5,672'

$ws.Range("A13").Value = 'Hotel Du Cadran'
$ws.Range("B13").Value = 'US$1,749'
$ws.Range("C13").Value = '8.5'
$ws.Range("D13").Value = 'Very Good'
$ws.Range("E13").Value = '795'

$ws.Range("A14").Value = 'Residhome Paris Rosa Parks'
$ws.Range("B14").Value = 'US$872'
$ws.Range("C14").Value = '8.2'
$ws.Range("E14").Value = '2,762'

$ws.Range("A15").Value = 'Hotel Anya'
$ws.Range("B15").Value = 'US$822'
$ws.Range("C15").Value = '7.4'
$ws.Range("D15").Value = 'Good'
$ws.Range("E15").Value = '1,132'

$ws.Range("A16").Value = 'Zoku Paris'
$ws.Range("B16").Value = 'US$1,806'
$ws.Range("C16").Value = '9.1'
$ws.Range("D16").Value = 'Wonderful'
$ws.Range("E16").Value = '884'

$ws.Range("A17").Value = 'Hotel Agenor'
$ws.Range("B17").Value = 'US$968'
$ws.Range("C17").Value = '7.0'
$ws.Range("D17").Value = 'Good'
$ws.Range("E17").Value = '1,242'

$ws.Range("A18").Value = 'Hôtel De Castiglione'
$ws.Range("B18").Value = 'US$1,833'
$ws.Range("C18").Value = '7.5'
$ws.Range("D18").Value = 'Good'
$ws.Range("E18").Value = '3,691'

$ws.Range("A19").Value = 'Monceau Elysées'
$ws.Range("B19").Value = 'US$1,686'
$ws.Range("C19").Value = '8.5'
$ws.Range("D19").Value = 'Very Good'
$ws.Range("E19").Value = '1,179'

$ws.Range("A20").Value = 'City Inn Paris'
$ws.Range("B20").Value = 'US$281'
$ws.Range("C20").Value = '6.7'
$ws.Range("D20").Value = 'Review score'
$ws.Range("E20").Value = '2,496'

$ws.Range("A21").Value = 'Le Coffice Auberge de Jeunesse'
$ws.Range("B21").Value = 'US$273'
$ws.Range("C21").Value = '5.3'
$ws.Range("D21").Value = 'Review score'
$ws.Range("E21").Value = '2,381'

$ws.Range("A22").Value = 'Austin''s Arts Et Metiers Hotel'
$ws.Range("B22").Value = 'US$1,653'
$ws.Range("C22").Value = '8.2'
$ws.Range("D22").Value = 'Very Good'
$ws.Range("E22").Value = '2,025'

$ws.Range("A23").Value = 'ibis Styles Paris Buttes Chaumont'
$ws.Range("B23").Value = 'US$1,014'
$ws.Range("C23").Value = '8.0'
$ws.Range("D23").Value = 'Very Good'
$ws.Range("E23").Value = '1,320'

$ws.Range("A24").Value = 'Glasgow Monceau by Patrick Hayat'
$ws.Range("B24").Value = 'US$1,166'
$ws.Range("C24").Value = '8.0'
$ws.Range("E24").Value = '674'

$ws.Range("A25").Value = 'Hotel Des Arts'
$ws.Range("B25").Value = 'US$576'
$ws.Range("C25").Value = '4.2'
$ws.Range("D25").Value = 'Review score'
$ws.Range("E25").Value = '1,049'

$ws.Range("A26").Value = 'Timhotel Montmartre'
$ws.Range("B26").Value = 'US$2,361'
$ws.Range("C26").Value = '8.3'
$ws.Range("D26").Value = 'Very Good'
$ws.Range("E26").Value = '2,282'

# Row 27 (Europe Saint Severin-Paris Notre Dame) was removed entirely,
# shrinking the used range from A1:E27 to A1:E26.
$ws.Rows(27).Delete()
